$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$notice = "填写须知：`n1、不能增加、删除列；`n2、不能修改灰色单元格；`n3、红色字段为必填字段，黑色字段为选填字段；`n4、删除餐厅编码后再导入，系统会根据ID删除没有餐厅编码的记录；`n5、删除餐厅编码后再导入，系统会根据ID删除没有餐厅编码的记录；`n6、删除餐厅编码后再导入，系统会根据ID删除没有餐厅编码的记录；删除餐厅编码后再导入，系统会根据ID删除没有餐厅编码的记录删除餐厅编码后再导入，系统会根据ID删除没有餐厅编码的记录`n7、删除餐厅编码后再导入，系统会根据ID删除没有餐厅编码的记录；"

$ws.Range("A1").Value = $notice
$ws.Rows.Item(1).RowHeight = 192.0
